# Inserts a new record row at row 441 (shifting the existing rows 441-471
# down to 442-472) in the "Hortaliza, Femacal de La Calera - Ají" sheet,
# then populates the newly inserted row with the new price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 441, pushing every row
# below it (441..471) down by one (442..472).
$ws.Rows.Item(441).Insert()

# Fill in the new row 441 with the new data record.
$ws.Cells.Item(441, 1).Value  = 3
$ws.Cells.Item(441, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(441, 3).Value  = "Coquimbo"
$ws.Cells.Item(441, 4).Value  = 44714
$ws.Cells.Item(441, 5).Value  = 5
$ws.Cells.Item(441, 6).Value  = 100112021
$ws.Cells.Item(441, 7).Value  = "Ají"
$ws.Cells.Item(441, 8).Value  = "Inferno"
$ws.Cells.Item(441, 9).Value  = "Primera"
$ws.Cells.Item(441, 10).Value = 76
$ws.Cells.Item(441, 11).Value = 22000
$ws.Cells.Item(441, 12).Value = 23000
$ws.Cells.Item(441, 13).Value = 22500
$ws.Cells.Item(441, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(441, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(441, 16).Value = 1500
$ws.Cells.Item(441, 17).Value = 15
$ws.Cells.Item(441, 18).Value = "Hortaliza"
